# Fill in the "Definition" column (D) for the concept rows on the
# "Concepts" sheet with the same text already present in the "Display"
# column (C), for rows 2 through 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $display = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value = $display
}
